# Update crypto "Price" (D) and "Volume(1h)" (E) columns with refreshed
# quote data. Values are stored as text in the workbook (e.g. "301.68",
# "0.56%"), so each assignment is given a leading apostrophe to force
# Excel to keep it as literal text instead of auto-converting it to a
# number or percentage.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'301.68"
$ws.Range("E2").Value = "'0.56%"
$ws.Range("D3").Value = "'32.14"
$ws.Range("E3").Value = "'1.33%"
$ws.Range("E4").Value = "'-2.93%"
$ws.Range("D5").Value = "'0.07898"
$ws.Range("E5").Value = "'-2.52%"
$ws.Range("D6").Value = "'2.097"
$ws.Range("E6").Value = "'-20.01%"
$ws.Range("D7").Value = "'7.796"
$ws.Range("E7").Value = "'0.17%"
$ws.Range("D8").Value = "'3.822"
$ws.Range("E8").Value = "'-2.17%"
$ws.Range("D9").Value = "'0.9290"
$ws.Range("E9").Value = "'-0.26%"
$ws.Range("D10").Value = "'0.1752"
$ws.Range("E10").Value = "'-0.53%"
$ws.Range("D11").Value = "'0.07983"
$ws.Range("E11").Value = "'8.36%"
$ws.Range("D12").Value = "'0.08573"
$ws.Range("E12").Value = "'-3.04%"
$ws.Range("D13").Value = "'0.03154"
$ws.Range("E13").Value = "'4.14%"
$ws.Range("D14").Value = "'0.1002"
$ws.Range("E14").Value = "'0.28%"
$ws.Range("D15").Value = "'0.001514"
$ws.Range("E15").Value = "'0.25%"
$ws.Range("D16").Value = "'0.005825"
$ws.Range("E16").Value = "'0.32%"
$ws.Range("E17").Value = "'2,100.33%"
$ws.Range("D18").Value = "'3.464"
$ws.Range("E18").Value = "'-2.97%"
$ws.Range("E19").Value = "'-0.34%"
$ws.Range("E20").Value = "'0.49%"
$ws.Range("E21").Value = "'-2.27%"
$ws.Range("E22").Value = "'2.81%"
$ws.Range("D24").Value = "'0.04600"
$ws.Range("E24").Value = "'-0.33%"
$ws.Range("E25").Value = "'-0.02%"
$ws.Range("D26").Value = "'0.004456"
$ws.Range("E26").Value = "'-1.60%"
$ws.Range("E27").Value = "'4.26%"
$ws.Range("D39").Value = "'0.01713"
$ws.Range("E39").Value = "'-2.69%"
$ws.Range("D40").Value = "'0.04771"
$ws.Range("E40").Value = "'3.72%"
$ws.Range("D41").Value = "'0.007448"
$ws.Range("E41").Value = "'8.02%"
$ws.Range("D42").Value = "'0.1359"
$ws.Range("E42").Value = "'-1.16%"
$ws.Range("D43").Value = "'0.002321"
$ws.Range("E43").Value = "'6.03%"
$ws.Range("D44").Value = "'0.01024"
$ws.Range("E44").Value = "'-0.60%"
$ws.Range("D45").Value = "'0.00005999"
$ws.Range("E45").Value = "'-5.12%"
$ws.Range("E46").Value = "'0.09%"
$ws.Range("D47").Value = "'0.003392"
$ws.Range("E47").Value = "'-59.60%"
$ws.Range("D48").Value = "'0.8234"
$ws.Range("E48").Value = "'9.99%"
$ws.Range("E49").Value = "'0.09%"
$ws.Range("E50").Value = "'0.09%"
